# Updated cryptos list (refreshed Price / Volume(1h) figures, plus the
# PEPE / InternetComputer(DFINITY) row swap at 26/27).
#
# Note: a few Price values (D9, D12, D19, D26, D28, D44, D45, D50) are
# numeric-looking strings that must keep their exact textual formatting
# (trailing zeros / significant digits), e.g. "0.500" or "1.00". Assigning
# them as a plain numeric-looking string would make Excel coerce the cell
# to a number and silently normalize it (e.g. "0.500" -> 0.5). Prefixing
# with a leading apostrophe forces Excel to keep them as text, matching
# the original author's intent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.304.98'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '3.690.36'
$ws.Range('E3').Value = '  -2.76%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '684.65'
$ws.Range('E5').Value = '  -2.90%  '
$ws.Range('E6').Value = '  -4.41%  '
$ws.Range('D7').Value = '3.688.52'
$ws.Range('E7').Value = '  -2.78%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.500'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('E10').Value = '  -7.16%  '
$ws.Range('D11').Value = '7.23'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('D12').Value = '''0.450'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('E13').Value = '  -6.87%  '
$ws.Range('D14').Value = '33.51'
$ws.Range('E14').Value = '  -7.24%  '
$ws.Range('D15').Value = '4.310.45'
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').Value = '3.695.29'
$ws.Range('E16').Value = '  -3.34%  '
$ws.Range('D17').Value = '69.379.20'
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D19').Value = '''16.30'
$ws.Range('E19').Value = '  -6.24%  '
$ws.Range('D20').Value = '6.62'
$ws.Range('E20').Value = '  -7.17%  '
$ws.Range('D21').Value = '482.25'
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').Value = '9.77'
$ws.Range('E22').Value = '  -8.09%  '
$ws.Range('D23').Value = '0.665'
$ws.Range('E23').Value = '  -8.65%  '
$ws.Range('D24').Value = '79.98'
$ws.Range('E24').Value = '  -5.90%  '
$ws.Range('D25').Value = '3.834.68'
$ws.Range('E25').Value = '  -2.89%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000128'
$ws.Range('E26').Value = '  -11.29%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '11.51'
$ws.Range('E27').Value = '  -4.69%  '
$ws.Range('D28').Value = '''1.00'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '9.57'
$ws.Range('E29').Value = '  -8.54%  '
$ws.Range('D30').Value = '1.84'
$ws.Range('E30').Value = '  -10.25%  '
$ws.Range('D31').Value = '2.76'
$ws.Range('E31').Value = '  -10.68%  '
$ws.Range('E32').Value = '  -5.01%  '
$ws.Range('E33').Value = '  -7.43%  '
$ws.Range('D34').Value = '26.98'
$ws.Range('E34').Value = '  -7.27%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = '0.165'
$ws.Range('E36').Value = '  -5.30%  '
$ws.Range('D37').Value = '3.650.73'
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('D38').Value = '8.56'
$ws.Range('E38').Value = '  -5.64%  '
$ws.Range('D39').Value = '6.05'
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').Value = '0.0941'
$ws.Range('E40').Value = '  -7.02%  '
$ws.Range('E42').Value = '  -6.01%  '
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').Value = '''0.960'
$ws.Range('E44').Value = '  -7.56%  '
$ws.Range('D45').Value = '''157.30'
$ws.Range('E45').Value = '  -4.42%  '
$ws.Range('D46').Value = '48.17'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  -14.47%  '
$ws.Range('D48').Value = '393.21'
$ws.Range('E48').Value = '  -7.20%  '
$ws.Range('D49').Value = '0.000278'
$ws.Range('E49').Value = '  -12.61%  '
$ws.Range('D50').Value = '''1.30'
$ws.Range('E50').Value = '  -4.68%  '
$ws.Range('D51').Value = '8.07'
$ws.Range('E51').Value = '  -6.29%  '
